$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.063.63'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '2.598.83'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.62'
$ws.Range("E5").Value = '  +2.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.93'
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.597'
$ws.Range("E7").Value = '  -0.61%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.90'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.41'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0837'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.11'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").Value = '2.996.92'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '2.608.11'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.913'
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.82'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '46.137.45'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").Value = '  -4.39%  '
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '279.26'
$ws.Range("E23").Value = '  +9.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.03'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.88'
$ws.Range("E27").Value = '  +6.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.06'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.71'
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  -3.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.84'
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("E33").Value = '  +1.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.58'
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.87'
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0838'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("E37").Value = '  -5.37%  '
$ws.Range("E38").Value = '  -4.41%  '
$ws.Range("E39").Value = '  +5.38%  '
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.00'
$ws.Range("E41").Value = '  +24.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.70'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0330'
$ws.Range("E43").Value = '  +2.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.58'
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.96'
$ws.Range("E45").Value = '  -5.69%  '
$ws.Range("D46").Value = '2.096.05'
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '95.58'
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("E49").Value = '  +5.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '108.60'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.200'
$ws.Range("E51").Value = '  -0.48%  '
